$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 45588.832
$ws.Range("J93").Value = 45588.832
$ws.Range("L93").Value = 45588.832
$ws.Range("N93").Value = -50580.832
$ws.Range("H95").Value = 37984
$ws.Range("J95").Value = 37984
$ws.Range("L95").Value = 37984
$ws.Range("N95").Value = -43476
$ws.Range("H123").Value = 37250
$ws.Range("J123").Value = 37250
$ws.Range("L123").Value = 37250
$ws.Range("N123").Value = -47050

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2695.17
$ws.Range("I32").Value = 1939.7472
$ws.Range("K32").Value = 1939.7472
$ws.Range("M32").Value = -1652.7472
$ws.Range("H106").Value = 46496
$ws.Range("J106").Value = 46496
$ws.Range("L106").Value = 46496
$ws.Range("N106").Value = -49020
$ws.Range("H107").Value = 37738
$ws.Range("J107").Value = 37738
$ws.Range("L107").Value = 37738
$ws.Range("N107").Value = -45418
$ws.Range("H120").Value = 39530.668
$ws.Range("J120").Value = 39530.668
$ws.Range("L120").Value = 39530.668
$ws.Range("N120").Value = -49206.668
$ws.Range("H121").Value = 45241
$ws.Range("J121").Value = 45241
$ws.Range("L121").Value = 45241
$ws.Range("N121").Value = -48735

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 45996
$ws.Range("J112").Value = 45996
$ws.Range("L112").Value = 45996
$ws.Range("N112").Value = -48950

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 184691.58
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 184691.58
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 184691.58
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -185281.58
$ws.Range("H34").Value = 184691.58
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 184691.58
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 184691.58
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -185095.58
$ws.Range("H100").Value = 47315.75
$ws.Range("J100").Value = 47315.75
$ws.Range("L100").Value = 47315.75
$ws.Range("N100").Value = -49479.75
$ws.Range("H104").Value = 28396.3
$ws.Range("J104").Value = 28396.3
$ws.Range("L104").Value = 28396.3
$ws.Range("N104").Value = -33638.3
$ws.Range("H110").Value = 41096.8
$ws.Range("J110").Value = 41096.8
$ws.Range("L110").Value = 41096.8
$ws.Range("N110").Value = -49276.8
$ws.Range("H114").Value = 20000
$ws.Range("J114").Value = 20000
$ws.Range("L114").Value = 20000
$ws.Range("N114").Value = -28678
$ws.Range("H115").Value = 28238.25
$ws.Range("J115").Value = 28238.25
$ws.Range("L115").Value = 28238.25
$ws.Range("N115").Value = -30588.25
$ws.Range("H119").Value = 44686
$ws.Range("J119").Value = 44686
$ws.Range("L119").Value = 44686
$ws.Range("N119").Value = -54362
$ws.Range("H124").Value = 45326
$ws.Range("J124").Value = 45326
$ws.Range("L124").Value = 45326
$ws.Range("N124").Value = -50236
$ws.Range("H125").Value = 38663
$ws.Range("J125").Value = 38663
$ws.Range("L125").Value = 38663
$ws.Range("N125").Value = -43583
$ws.Range("H131").Value = 38318
$ws.Range("J131").Value = 38318
$ws.Range("L131").Value = 38318
$ws.Range("N131").Value = -48398
$ws.Range("H133").Value = 27999.62
$ws.Range("J133").Value = 27999.62
$ws.Range("L133").Value = 27999.62
$ws.Range("N133").Value = -33059.62

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 37261
$ws.Range("J96").Value = 37261
$ws.Range("L96").Value = 37261
$ws.Range("N96").Value = -42753
$ws.Range("H104").Value = 44542
$ws.Range("J104").Value = 44542
$ws.Range("L104").Value = 44542
$ws.Range("N104").Value = -51530
$ws.Range("H116").Value = 38998
$ws.Range("J116").Value = 38998
$ws.Range("L116").Value = 38998
$ws.Range("N116").Value = -48176
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 37992
$ws.Range("J120").Value = 37992
$ws.Range("L120").Value = 37992
$ws.Range("N120").Value = -47668
$ws.Range("H127").Value = 47300.332
$ws.Range("J127").Value = 47300.332
$ws.Range("L127").Value = 47300.332
$ws.Range("N127").Value = -57220.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 14500
$ws.Range("I74").Value = 14500
$ws.Range("K74").Value = 14500
$ws.Range("M74").Value = -13502
$ws.Range("H77").Value = 14500
$ws.Range("I77").Value = 14500
$ws.Range("K77").Value = 43500
$ws.Range("M77").Value = -38508
$ws.Range("H108").Value = 48584
$ws.Range("J108").Value = 48584
$ws.Range("L108").Value = 48584
$ws.Range("N108").Value = -56264
$ws.Range("H109").Value = 35277
$ws.Range("J109").Value = 35277
$ws.Range("L109").Value = 35277
$ws.Range("N109").Value = -38051
$ws.Range("H117").Value = 36892
$ws.Range("J117").Value = 36892
$ws.Range("L117").Value = 36892
$ws.Range("N117").Value = -46070
$ws.Range("H123").Value = 26662.5
$ws.Range("J123").Value = 32883.332
$ws.Range("L123").Value = 32883.332
$ws.Range("N123").Value = -42683.332
$ws.Range("H129").Value = 44421
$ws.Range("J129").Value = 44421
$ws.Range("L129").Value = 44421
$ws.Range("N129").Value = -54421
$ws.Range("H131").Value = 43318
$ws.Range("J131").Value = 43318
$ws.Range("L131").Value = 43318
$ws.Range("N131").Value = -53398

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45890.668
$ws.Range("J16").Value = 45890.668
$ws.Range("L16").Value = 45890.668
$ws.Range("N16").Value = -46474.668
$ws.Range("H27").Value = 22996
$ws.Range("J27").Value = 30492
$ws.Range("L27").Value = 30492
$ws.Range("N27").Value = -30630
$ws.Range("H109").Value = 34873
$ws.Range("J109").Value = 34873
$ws.Range("L109").Value = 34873
$ws.Range("N109").Value = -37647
$ws.Range("H115").Value = 37998
$ws.Range("J115").Value = 37998
$ws.Range("L115").Value = 37998
$ws.Range("N115").Value = -41132
$ws.Range("H117").Value = 47156
$ws.Range("J117").Value = 47156
$ws.Range("L117").Value = 47156
$ws.Range("N117").Value = -56334
$ws.Range("H118").Value = 38997.332
$ws.Range("J118").Value = 41997
$ws.Range("L118").Value = 41997
$ws.Range("N118").Value = -45311
$ws.Range("H119").Value = 46227.332
$ws.Range("J119").Value = 46227.332
$ws.Range("L119").Value = 46227.332
$ws.Range("N119").Value = -55903.332
$ws.Range("H120").Value = 45416
$ws.Range("J120").Value = 45416
$ws.Range("L120").Value = 45416
$ws.Range("N120").Value = -55092
$ws.Range("H127").Value = 16283.429
$ws.Range("J127").Value = 41984
$ws.Range("L127").Value = 41984
$ws.Range("N127").Value = -51904
$ws.Range("H129").Value = 32996
$ws.Range("J129").Value = 32996
$ws.Range("L129").Value = 32996
$ws.Range("N129").Value = -42996
